# Helper: VBA-style RGB() -> packed BGR long used by PowerPoint's RGB color objects.
function RGBColor($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation

# --- 1. Table on slide 5: switch its table style to the built-in style GUID ---
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $shp = $s5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{9E7974B8-AE52-4610-9771-7008DB89473B}")
    }
}

# --- 2. Slide master theme: recolor from "Integral" (Red Violet) to "Office Theme" ---
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
$tcs.Colors(1).RGB  = RGBColor 0x00 0x00 0x00   # dk1      -> 000000
$tcs.Colors(2).RGB  = RGBColor 0xFF 0xFF 0xFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = RGBColor 0x44 0x54 0x6A   # dk2      -> 44546A
$tcs.Colors(4).RGB  = RGBColor 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = RGBColor 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = RGBColor 0xED 0x7D 0x31   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = RGBColor 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = RGBColor 0xFF 0xC0 0x00   # accent4  -> FFC000
$tcs.Colors(9).RGB  = RGBColor 0x44 0x72 0xC4   # accent5  -> 4472C4
$tcs.Colors(10).RGB = RGBColor 0x70 0xAD 0x47   # accent6  -> 70AD47
$tcs.Colors(11).RGB = RGBColor 0x05 0x63 0xC1   # hlink    -> 0563C1
$tcs.Colors(12).RGB = RGBColor 0x95 0x4F 0x72   # folHlink -> 954F72
